{"js": "// Auto-generated: insert the new introduction paragraphs before the\n// existing (sole) paragraph in the document body.\nconst body = context.document.body;\n\nconst p0 = body.insertParagraph(\"Apache solar\", Word.InsertLocation.start);\n\nconst p1 = p0.insertParagraph(\"Introduction\", Word.InsertLocation.after);\n\nconst p2 = p1.insertParagraph(\"Solr was created in 2004 as an in house project to add search capability for the CNET networks. It was donated to APACHE Software Foundation after CNET Networks decided to make it open source. After a number of versions of the Solr, with every version having enhancements, the latest version today is the Solr 5.0 which is a standalone application. Solr is an open source \", Word.InsertLocation.after);\np2.insertText(\"search platform written in JAVA. It has full text search capabilities that enable matching phrases, hit highlighting, real time indexing\", Word.InsertLocation.end);\np2.insertText(\" so that the user can see the content whenever he wants to see\", Word.InsertLocation.end);\np2.insertText(\", \", Word.InsertLocation.end);\np2.insertText(\" no SQL features, extensible plugin architecture, dynamic clustering, database integration. Providing distributed search and index replication, Solr is highly scalable and fault tolerant. \", Word.InsertLocation.end);\n\nconst p3 = p2.insertParagraph(\"Why Apache Solr is so Popular?\", Word.InsertLocation.after);\n\nconst p4 = p3.insertParagraph(\"Although there are databases and frameworks such as HADOOP, Apache Solr has been on top of these due to its standalone noSQL store. Solr effectively serves as a data \", Word.InsertLocation.after);\np4.insertText(\"access layer for doing key value lookups as well as making the data fully indexed and searchable. \", Word.InsertLocation.end);\np4.insertText(\" Search is always the first requirement in many deployments, which is pushing many organisations to use the search engine like a noSQL store for that robustness and flexibility. \", Word.InsertLocation.end);\n\nconst p5 = p4.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add the introductory paper content as new paragraphs inserted right\n# before the document's existing (sole) paragraph, which holds the\n# '_GoBack' bookmark and must remain the last paragraph in the body.\n$d = $word.ActiveDocument\n\n# Helper pattern: the original paragraph is always the LAST paragraph in\n# the body, so re-querying Paragraphs(Count) each time and calling\n# InsertParagraphBefore() on it inserts a fresh blank paragraph directly\n# ahead of it -- i.e. right after whatever we inserted previously. That\n# lets us fill the new paragraphs in natural top-to-bottom order.\n\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphBefore()\n$r = $d.Paragraphs($d.Paragraphs.Count - 1).Range\n$r.InsertAfter('Apache solar')\n\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphBefore()\n$r = $d.Paragraphs($d.Paragraphs.Count - 1).Range\n$r.InsertAfter('Introduction')\n\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphBefore()\n$r = $d.Paragraphs($d.Paragraphs.Count - 1).Range\n$r.InsertAfter('Solr was created in 2004 as an in house project to add search capability for the CNET networks. It was donated to APACHE Software Foundation after CNET Networks decided to make it open source. After a number of versions of the Solr, with every version having enhancements, the latest version today is the Solr 5.0 which is a standalone application. Solr is an open source ')\n$r.InsertAfter('search platform written in JAVA. It has full text search capabilities that enable matching phrases, hit highlighting, real time indexing')\n$r.InsertAfter(' so that the user can see the content whenever he wants to see')\n$r.InsertAfter(', ')\n$r.InsertAfter(' no SQL features, extensible plugin architecture, dynamic clustering, database integration. Providing distributed search and index replication, Solr is highly scalable and fault tolerant. ')\n\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphBefore()\n$r = $d.Paragraphs($d.Paragraphs.Count - 1).Range\n$r.InsertAfter('Why Apache Solr is so Popular?')\n\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphBefore()\n$r = $d.Paragraphs($d.Paragraphs.Count - 1).Range\n$r.InsertAfter('Although there are databases and frameworks such as HADOOP, Apache Solr has been on top of these due to its standalone noSQL store. Solr effectively serves as a data ')\n$r.InsertAfter('access layer for doing key value lookups as well as making the data fully indexed and searchable. ')\n$r.InsertAfter(' Search is always the first requirement in many deployments, which is pushing many organisations to use the search engine like a noSQL store for that robustness and flexibility. ')\n\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphBefore()\n"}
